# "implemented new partition method - link process method"
# Update the per-balancing-authority solid/liquid/gas emission-factor
# numbers on Sheet1, and clear out the rows that no longer carry data
# (they keep their eGRID label in column A but lose the B:D figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- updated figures ---------------------------------------------------
$ws.Range("B5").Value  = 581.71667763057258
$ws.Range("C5").Value  = 9.37152389365842
$ws.Range("D5").Value  = 71.336529653564654

$ws.Range("B7").Value  = 1481.2459527721019
$ws.Range("C7").Value  = 8.9355307679654778
$ws.Range("D7").Value  = 216.14615524129863

$ws.Range("B8").Value  = 926.96648857857485
$ws.Range("C8").Value  = 11.915308011847486
$ws.Range("D8").Value  = 166.5283616205796

$ws.Range("B11").Value = 360.07685385388277
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 209.147062138102

$ws.Range("B12").Value = 426.89615119692996
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 181.0009148305198

$ws.Range("B14").Value = 524.03691033842938
$ws.Range("C14").Value = 6.4953839331559511
$ws.Range("D14").Value = 49.37064404925345

$ws.Range("B17").Value = 521.63422561664117
$ws.Range("C17").Value = 3.9401938910513108
$ws.Range("D17").Value = 168.46487077844267

$ws.Range("B18").Value = 719.34697608887632
$ws.Range("C18").Value = 3.8481375907931366
$ws.Range("D18").Value = 110.33988468012224

$ws.Range("B19").Value = 386.52174845250829
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 241.98307329129304

$ws.Range("B20").Value = 956.37655799275387
$ws.Range("C20").Value = 13.891290815888048
$ws.Range("D20").Value = 71.196545953163593

$ws.Range("B21").Value = 400.449468621535
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 224.05059710106119

$ws.Range("B22").Value = 399.5435836117606
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 215.04235918102626

$ws.Range("B23").Value = 385.8961752070943
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 248.04596408008157

$ws.Range("B24").Value = 370.25238895331245
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 259.01993733160361

$ws.Range("B25").Value = 378.78789099852276
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 216.75838487344981

$ws.Range("B26").Value = 673.85843676530908
$ws.Range("C26").Value = 12.533757321587807
$ws.Range("D26").Value = 137.10729282684588

$ws.Range("B27").Value = 1049.6795339494079
$ws.Range("C27").Value = 13.828836209725528
$ws.Range("D27").Value = 80.478181970815399

$ws.Range("B28").Value = 1266.2631935861396
$ws.Range("C28").Value = 14.381239009481749
$ws.Range("D28").Value = 110.06301747337619

# --- rows whose B:D figures are no longer populated --------------------
$ws.Range("B2:D2").ClearContents()
$ws.Range("B3:D3").ClearContents()
$ws.Range("B4:D4").ClearContents()
$ws.Range("B6:D6").ClearContents()
$ws.Range("B9:D9").ClearContents()
$ws.Range("B10:D10").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("B15:D15").ClearContents()
$ws.Range("B16:D16").ClearContents()

# --- reset the view: scroll back to the top, select the home cell ------
$ws.Range("A1").Select()
